$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.956.86"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.658.78"
$ws.Range("E3").Value = "  +2.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.35"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3891"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.17"
$ws.Range("E9").Value = "  +5.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.356"
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.0000"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08472"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.96"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.154"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.889"
$ws.Range("E15").Value = "  +6.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001308"
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.651.45"
$ws.Range("E17").Value = "  +2.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.83"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07003"
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.81"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.912"
$ws.Range("E21").Value = "  +1.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.61"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.933.56"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.488"
$ws.Range("E25").Value = "  +3.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.052"
$ws.Range("E26").Value = "  +8.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.08"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.50"
$ws.Range("E28").Value = "  -2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.435"
$ws.Range("E29").Value = "  +3.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "139.25"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.781"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.494"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.832.29"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.030"
$ws.Range("E34").Value = "  +8.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08053"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02955"
$ws.Range("E36").Value = "  +3.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.99"
$ws.Range("E37").Value = "  +5.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.668"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2686"
$ws.Range("E39").Value = "  +2.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09121"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7541"
$ws.Range("E41").Value = "  +1.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.47"
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.418"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.17"
$ws.Range("E44").Value = "  +2.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6949"
$ws.Range("E45").Value = "  +2.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.465"
$ws.Range("E46").Value = "  +1.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.077"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08286"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "134.14"
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("E51").Value = "  +7.33%  "
